# New crime data collected — weekly refresh of the 34th Precinct CompStat sheet.
# Updates the report header (volume number + covered week dates) and refreshes
# every weekly/28-day/YTD/2-year crime-complaint figure in rows 14-30, including
# a handful of cells that flip between a numeric count and the "0"/"***.*"
# placeholder text used when a category has no data for the period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: "Volume 30   Number  39" -> "...Number  40"
#         "Report Covering the Week  9/25/2023  Through  10/1/2023"
#      -> "...10/2/2023  Through  10/8/2023"
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "40"
$ws.Range("C9").Characters(27, 9).Text = "10/2/2023"
$ws.Range("C9").Characters(47, 9).Text = "10/8/2023"

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds a number into the shared
# "placeholder" text (e.g. "0" or "***.*"), reusing the number format/style
# of a donor cell that already displays that placeholder so the style index
# the workbook ends up with matches (right-aligned "General" text style).
# ---------------------------------------------------------------------------
function Set-PlaceholderText($target, $text, $donorRef) {
    $ws.Range($target).Value = "'" + $text
    $ws.Range($donorRef).Copy()
    $ws.Range($target).PasteSpecial(-4122)
}

# Helper: convert a cell that currently holds placeholder text back into a
# plain number, reusing the numeric style of a donor cell.
function Set-NumericValue($target, $value, $donorRef) {
    $ws.Range($target).Value = $value
    $ws.Range($donorRef).Copy()
    $ws.Range($target).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Row 14
# ---------------------------------------------------------------------------
Set-PlaceholderText "C14" "0" "D14"

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 21
$ws.Range("K15").Value = -66.666666666666
$ws.Range("N15").Value = -90.277777777777

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -57.142857142857
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = -58.064516129032
$ws.Range("I16").Value = 166
$ws.Range("J16").Value = 218
$ws.Range("K16").Value = -23.853211009174
$ws.Range("L16").Value = 6.410256410256
$ws.Range("M16").Value = -20.192307692307
$ws.Range("N16").Value = -85.309734513274

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 22.727272727272
$ws.Range("I17").Value = 256
$ws.Range("J17").Value = 247
$ws.Range("K17").Value = 3.643724696356
$ws.Range("L17").Value = 13.777777777777
$ws.Range("M17").Value = 52.380952380952
$ws.Range("N17").Value = -65.591397849462

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
Set-PlaceholderText "C18" "0" "D14"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -45.454545454545
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = -0.862068965517
$ws.Range("L18").Value = 38.554216867469
$ws.Range("M18").Value = -17.266187050359
$ws.Range("N18").Value = -93.072289156626

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -8.510638297872
$ws.Range("I19").Value = 415
$ws.Range("J19").Value = 456
$ws.Range("K19").Value = -8.991228070175
$ws.Range("L19").Value = -9.586056644880
$ws.Range("M19").Value = 52.014652014652
$ws.Range("N19").Value = -58.582834331337

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -42.857142857142
$ws.Range("F20").Value = 13
$ws.Range("H20").Value = 8.333333333333
$ws.Range("I20").Value = 178
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = -28.8
$ws.Range("L20").Value = 7.878787878787
$ws.Range("M20").Value = 165.671641791045
$ws.Range("N20").Value = -88.388780169602

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -40.540540540540
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 126
$ws.Range("H21").Value = -17.460317460317
$ws.Range("I21").Value = 1140
$ws.Range("J21").Value = 1312
$ws.Range("K21").Value = -13.109756097561
$ws.Range("L21").Value = 2.517985611510
$ws.Range("M21").Value = 31.034482758620
$ws.Range("N21").Value = -81.586173477628

# ---------------------------------------------------------------------------
# Row 22 — C22 becomes the placeholder "0" while D22/E22 flip from
# placeholder text back to real numbers (the blank/placeholder slot moved
# from the D/E columns to C).
# ---------------------------------------------------------------------------
Set-PlaceholderText "C22" "0" "D14"
Set-NumericValue "D22" 2 "D18"
Set-NumericValue "E22" -100 "E15"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 47.619047619047

# ---------------------------------------------------------------------------
# Row 23 — same placeholder shuffle as row 22, but on both the C/D/E triplet
# and the F/G/H triplet.
# ---------------------------------------------------------------------------
Set-PlaceholderText "C23" "0" "D14"
Set-NumericValue "D23" 2 "D18"
Set-NumericValue "E23" -100 "E15"
$ws.Range("F23").Value = 4
Set-NumericValue "G23" 2 "D18"
Set-NumericValue "H23" 100 "E15"
$ws.Range("J23").Value = 32
$ws.Range("K23").Value = -18.75
$ws.Range("M23").Value = 62.5

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -47.826086956521
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -16.666666666666
$ws.Range("I24").Value = 1064
$ws.Range("J24").Value = 1109
$ws.Range("K24").Value = -4.057709648331
$ws.Range("L24").Value = 57.863501483679
$ws.Range("M24").Value = 121.205821205821

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 166.666666666667
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 45.945945945945
$ws.Range("I25").Value = 420
$ws.Range("J25").Value = 351
$ws.Range("K25").Value = 19.658119658119
$ws.Range("L25").Value = 28.048780487804
$ws.Range("M25").Value = -1.639344262295

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("J26").Value = 32
$ws.Range("K26").Value = -46.875

# ---------------------------------------------------------------------------
# Row 27 — C27/D27/E27 flip from placeholder text to real numbers.
# ---------------------------------------------------------------------------
Set-NumericValue "C27" 1 "D18"
Set-NumericValue "D27" 3 "D18"
Set-NumericValue "E27" -66.666666666666 "E15"
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = 25.714285714285
$ws.Range("L27").Value = -10.204081632653

# ---------------------------------------------------------------------------
# Row 28 — C28/D28/E28 flip from real numbers to placeholder text.
# ---------------------------------------------------------------------------
Set-PlaceholderText "C28" "0" "D14"
Set-PlaceholderText "D28" "0" "D14"
Set-PlaceholderText "E28" "***.*" "E14"
$ws.Range("M28").Value = 71.428571428571
$ws.Range("N28").Value = -92.356687898089

# ---------------------------------------------------------------------------
# Row 29 — same as row 28.
# ---------------------------------------------------------------------------
Set-PlaceholderText "C29" "0" "D14"
Set-PlaceholderText "D29" "0" "D14"
Set-PlaceholderText "E29" "***.*" "E14"
$ws.Range("M29").Value = 57.142857142857
$ws.Range("N29").Value = -91.970802919708

# ---------------------------------------------------------------------------
# Row 30 — G30/H30 flip from real numbers to placeholder text.
# ---------------------------------------------------------------------------
Set-PlaceholderText "G30" "0" "D14"
Set-PlaceholderText "H30" "***.*" "E14"
